$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 30 and 31 (columns F:V) ---
$row30 = @()
$row31 = @()
for ($col = 6; $col -le 22; $col++) {
    $row30 += ,($ws.Cells.Item(30, $col).Value())
    $row31 += ,($ws.Cells.Item(31, $col).Value())
}
for ($i = 0; $i -lt $row30.Length; $i++) {
    $col = 6 + $i
    $ws.Cells.Item(30, $col).Value = $row31[$i]
    $ws.Cells.Item(31, $col).Value = $row30[$i]
}

# --- Append new rows 60-77 (A = 59-76) ---
# Row 60
$ws.Cells.Item(60, 1).Value = 59
$ws.Cells.Item(60, 2).Value = "malta"
$ws.Cells.Item(60, 3).Value = "premier-league"
$ws.Cells.Item(60, 4).Value = "2023-2024"
$ws.Cells.Item(60, 5).Value = 45263.45833333334
$ws.Cells.Item(60, 6).Value = "Naxxar"
$ws.Cells.Item(60, 7).Value = 2
$ws.Cells.Item(60, 8).Value = "Gudja"
$ws.Cells.Item(60, 9).Value = 1
$ws.Cells.Item(60, 10).Value = 2.09
$ws.Cells.Item(60, 11).Value = "01/12/2023 23:12"
$ws.Cells.Item(60, 12).Value = 1.75
$ws.Cells.Item(60, 13).Value = "03/12/2023 10:59"
$ws.Cells.Item(60, 14).Value = 3.03
$ws.Cells.Item(60, 15).Value = "01/12/2023 23:12"
$ws.Cells.Item(60, 16).Value = 3.56
$ws.Cells.Item(60, 17).Value = "03/12/2023 10:59"
$ws.Cells.Item(60, 18).Value = 3.38
$ws.Cells.Item(60, 19).Value = "01/12/2023 23:12"
$ws.Cells.Item(60, 20).Value = 4.56
$ws.Cells.Item(60, 21).Value = "03/12/2023 10:59"
$ws.Cells.Item(60, 22).Value = "https://www.betexplorer.com/football/malta/premier-league/naxxar-lions-gudja/GC4rUKnk/"

# Row 61
$ws.Cells.Item(61, 1).Value = 60
$ws.Cells.Item(61, 2).Value = "malta"
$ws.Cells.Item(61, 3).Value = "premier-league"
$ws.Cells.Item(61, 4).Value = "2023-2024"
$ws.Cells.Item(61, 5).Value = 45263.58333333334
$ws.Cells.Item(61, 6).Value = "Floriana"
$ws.Cells.Item(61, 7).Value = 1
$ws.Cells.Item(61, 8).Value = "Valletta"
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 10).Value = 1.93
$ws.Cells.Item(61, 11).Value = "02/12/2023 02:13"
$ws.Cells.Item(61, 12).Value = 1.87
$ws.Cells.Item(61, 13).Value = "03/12/2023 13:07"
$ws.Cells.Item(61, 14).Value = 3.11
$ws.Cells.Item(61, 15).Value = "02/12/2023 02:13"
$ws.Cells.Item(61, 16).Value = 3.18
$ws.Cells.Item(61, 17).Value = "03/12/2023 13:07"
$ws.Cells.Item(61, 18).Value = 3.75
$ws.Cells.Item(61, 19).Value = "02/12/2023 02:13"
$ws.Cells.Item(61, 20).Value = 4.53
$ws.Cells.Item(61, 21).Value = "03/12/2023 13:07"
$ws.Cells.Item(61, 22).Value = "https://www.betexplorer.com/football/malta/premier-league/floriana-valletta/lQ2jSbH1/"

# Row 62
$ws.Cells.Item(62, 1).Value = 61
$ws.Cells.Item(62, 2).Value = "malta"
$ws.Cells.Item(62, 3).Value = "premier-league"
$ws.Cells.Item(62, 4).Value = "2023-2024"
$ws.Cells.Item(62, 5).Value = 45263.58333333334
$ws.Cells.Item(62, 6).Value = "Gzira"
$ws.Cells.Item(62, 7).Value = 2
$ws.Cells.Item(62, 8).Value = "Marsaxlokk"
$ws.Cells.Item(62, 9).Value = 4
$ws.Cells.Item(62, 10).Value = 2.05
$ws.Cells.Item(62, 11).Value = "02/12/2023 02:13"
$ws.Cells.Item(62, 12).Value = 2.42
$ws.Cells.Item(62, 13).Value = "03/12/2023 13:55"
$ws.Cells.Item(62, 14).Value = 3.12
$ws.Cells.Item(62, 15).Value = "02/12/2023 02:13"
$ws.Cells.Item(62, 16).Value = 2.66
$ws.Cells.Item(62, 17).Value = "03/12/2023 13:55"
$ws.Cells.Item(62, 18).Value = 3.36
$ws.Cells.Item(62, 19).Value = "02/12/2023 02:13"
$ws.Cells.Item(62, 20).Value = 3.54
$ws.Cells.Item(62, 21).Value = "03/12/2023 13:56"
$ws.Cells.Item(62, 22).Value = "https://www.betexplorer.com/football/malta/premier-league/gzira-marsaxlokk/f33nTv2e/"

# Row 63
$ws.Cells.Item(63, 1).Value = 62
$ws.Cells.Item(63, 2).Value = "malta"
$ws.Cells.Item(63, 3).Value = "premier-league"
$ws.Cells.Item(63, 4).Value = "2023-2024"
$ws.Cells.Item(63, 5).Value = 45263.67708333334
$ws.Cells.Item(63, 6).Value = "Hibernians"
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = "Sliema"
$ws.Cells.Item(63, 9).Value = 1
$ws.Cells.Item(63, 10).Value = 2.19
$ws.Cells.Item(63, 11).Value = "02/12/2023 04:43"
$ws.Cells.Item(63, 12).Value = 2.67
$ws.Cells.Item(63, 13).Value = "03/12/2023 16:06"
$ws.Cells.Item(63, 14).Value = 3.45
$ws.Cells.Item(63, 15).Value = "02/12/2023 04:43"
$ws.Cells.Item(63, 16).Value = 2.7
$ws.Cells.Item(63, 17).Value = "03/12/2023 16:06"
$ws.Cells.Item(63, 18).Value = 2.79
$ws.Cells.Item(63, 19).Value = "02/12/2023 04:43"
$ws.Cells.Item(63, 20).Value = 3.07
$ws.Cells.Item(63, 21).Value = "03/12/2023 16:06"
$ws.Cells.Item(63, 22).Value = "https://www.betexplorer.com/football/malta/premier-league/hibernians-sliema/006fRIW7/"

# Row 64
$ws.Cells.Item(64, 1).Value = 63
$ws.Cells.Item(64, 2).Value = "malta"
$ws.Cells.Item(64, 3).Value = "premier-league"
$ws.Cells.Item(64, 4).Value = "2023-2024"
$ws.Cells.Item(64, 5).Value = 45269.58333333334
$ws.Cells.Item(64, 6).Value = "Sliema"
$ws.Cells.Item(64, 7).Value = 1
$ws.Cells.Item(64, 8).Value = "Balzan"
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 2.02
$ws.Cells.Item(64, 11).Value = "08/12/2023 02:13"
$ws.Cells.Item(64, 12).Value = 2.19
$ws.Cells.Item(64, 13).Value = "09/12/2023 13:51"
$ws.Cells.Item(64, 14).Value = 3.16
$ws.Cells.Item(64, 15).Value = "08/12/2023 02:13"
$ws.Cells.Item(64, 16).Value = 2.77
$ws.Cells.Item(64, 17).Value = "09/12/2023 13:51"
$ws.Cells.Item(64, 18).Value = 3.4
$ws.Cells.Item(64, 19).Value = "08/12/2023 02:13"
$ws.Cells.Item(64, 20).Value = 3.95
$ws.Cells.Item(64, 21).Value = "09/12/2023 13:51"
$ws.Cells.Item(64, 22).Value = "https://www.betexplorer.com/football/malta/premier-league/sliema-balzan-fc/tdqskMOE/"

# Row 65
$ws.Cells.Item(65, 1).Value = 64
$ws.Cells.Item(65, 2).Value = "malta"
$ws.Cells.Item(65, 3).Value = "premier-league"
$ws.Cells.Item(65, 4).Value = "2023-2024"
$ws.Cells.Item(65, 5).Value = 45269.58333333334
$ws.Cells.Item(65, 6).Value = "Valletta"
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = "Hibernians"
$ws.Cells.Item(65, 9).Value = 2
$ws.Cells.Item(65, 10).Value = 2.68
$ws.Cells.Item(65, 11).Value = "08/12/2023 02:13"
$ws.Cells.Item(65, 12).Value = 3.11
$ws.Cells.Item(65, 13).Value = "09/12/2023 13:55"
$ws.Cells.Item(65, 14).Value = 2.88
$ws.Cells.Item(65, 15).Value = "08/12/2023 02:13"
$ws.Cells.Item(65, 16).Value = 2.63
$ws.Cells.Item(65, 17).Value = "09/12/2023 13:55"
$ws.Cells.Item(65, 18).Value = 2.6
$ws.Cells.Item(65, 19).Value = "08/12/2023 02:13"
$ws.Cells.Item(65, 20).Value = 2.7
$ws.Cells.Item(65, 21).Value = "09/12/2023 13:55"
$ws.Cells.Item(65, 22).Value = "https://www.betexplorer.com/football/malta/premier-league/valletta-hibernians/84CAXggs/"

# Row 66
$ws.Cells.Item(66, 1).Value = 65
$ws.Cells.Item(66, 2).Value = "malta"
$ws.Cells.Item(66, 3).Value = "premier-league"
$ws.Cells.Item(66, 4).Value = "2023-2024"
$ws.Cells.Item(66, 5).Value = 45269.70833333334
$ws.Cells.Item(66, 6).Value = "Mosta"
$ws.Cells.Item(66, 7).Value = 0
$ws.Cells.Item(66, 8).Value = "Floriana"
$ws.Cells.Item(66, 9).Value = 3
$ws.Cells.Item(66, 10).Value = 6.24
$ws.Cells.Item(66, 11).Value = "08/12/2023 05:12"
$ws.Cells.Item(66, 12).Value = 7.42
$ws.Cells.Item(66, 13).Value = "09/12/2023 16:35"
$ws.Cells.Item(66, 14).Value = 3.91
$ws.Cells.Item(66, 15).Value = "08/12/2023 05:12"
$ws.Cells.Item(66, 16).Value = 4.23
$ws.Cells.Item(66, 17).Value = "09/12/2023 16:35"
$ws.Cells.Item(66, 18).Value = 1.45
$ws.Cells.Item(66, 19).Value = "08/12/2023 05:12"
$ws.Cells.Item(66, 20).Value = 1.43
$ws.Cells.Item(66, 21).Value = "09/12/2023 16:35"
$ws.Cells.Item(66, 22).Value = "https://www.betexplorer.com/football/malta/premier-league/mosta-fc-floriana/neBEWD8m/"

# Row 67
$ws.Cells.Item(67, 1).Value = 66
$ws.Cells.Item(67, 2).Value = "malta"
$ws.Cells.Item(67, 3).Value = "premier-league"
$ws.Cells.Item(67, 4).Value = "2023-2024"
$ws.Cells.Item(67, 5).Value = 45270.45833333334
$ws.Cells.Item(67, 6).Value = "Marsaxlokk"
$ws.Cells.Item(67, 7).Value = 3
$ws.Cells.Item(67, 8).Value = "Santa Lucia"
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 10).Value = 1.51
$ws.Cells.Item(67, 11).Value = "08/12/2023 23:12"
$ws.Cells.Item(67, 12).Value = 1.62
$ws.Cells.Item(67, 13).Value = "10/12/2023 10:52"
$ws.Cells.Item(67, 14).Value = 3.75
$ws.Cells.Item(67, 15).Value = "08/12/2023 23:12"
$ws.Cells.Item(67, 16).Value = 3.32
$ws.Cells.Item(67, 17).Value = "10/12/2023 10:52"
$ws.Cells.Item(67, 18).Value = 5.23
$ws.Cells.Item(67, 19).Value = "08/12/2023 23:12"
$ws.Cells.Item(67, 20).Value = 6.55
$ws.Cells.Item(67, 21).Value = "10/12/2023 10:52"
$ws.Cells.Item(67, 22).Value = "https://www.betexplorer.com/football/malta/premier-league/marsaxlokk-santa-lucia/UJ0JVXNg/"

# Row 68
$ws.Cells.Item(68, 1).Value = 67
$ws.Cells.Item(68, 2).Value = "malta"
$ws.Cells.Item(68, 3).Value = "premier-league"
$ws.Cells.Item(68, 4).Value = "2023-2024"
$ws.Cells.Item(68, 5).Value = 45270.58333333334
$ws.Cells.Item(68, 6).Value = "Sirens"
$ws.Cells.Item(68, 7).Value = 1
$ws.Cells.Item(68, 8).Value = "Naxxar"
$ws.Cells.Item(68, 9).Value = 1
$ws.Cells.Item(68, 10).Value = 3.25
$ws.Cells.Item(68, 11).Value = "09/12/2023 02:13"
$ws.Cells.Item(68, 12).Value = 2.49
$ws.Cells.Item(68, 13).Value = "10/12/2023 13:58"
$ws.Cells.Item(68, 14).Value = 3.05
$ws.Cells.Item(68, 15).Value = "09/12/2023 02:13"
$ws.Cells.Item(68, 16).Value = 3.47
$ws.Cells.Item(68, 17).Value = "10/12/2023 13:55"
$ws.Cells.Item(68, 18).Value = 2.13
$ws.Cells.Item(68, 19).Value = "09/12/2023 02:13"
$ws.Cells.Item(68, 20).Value = 2.62
$ws.Cells.Item(68, 21).Value = "10/12/2023 13:58"
$ws.Cells.Item(68, 22).Value = "https://www.betexplorer.com/football/malta/premier-league/sirens-naxxar-lions/rN4NUiw0/"

# Row 69
$ws.Cells.Item(69, 1).Value = 68
$ws.Cells.Item(69, 2).Value = "malta"
$ws.Cells.Item(69, 3).Value = "premier-league"
$ws.Cells.Item(69, 4).Value = "2023-2024"
$ws.Cells.Item(69, 5).Value = 45273.45833333334
$ws.Cells.Item(69, 6).Value = "Gudja"
$ws.Cells.Item(69, 7).Value = 1
$ws.Cells.Item(69, 8).Value = "Birkirkara"
$ws.Cells.Item(69, 9).Value = 1
$ws.Cells.Item(69, 10).Value = 6.2
$ws.Cells.Item(69, 11).Value = "11/12/2023 23:12"
$ws.Cells.Item(69, 12).Value = 5.56
$ws.Cells.Item(69, 13).Value = "13/12/2023 10:55"
$ws.Cells.Item(69, 14).Value = 3.8
$ws.Cells.Item(69, 15).Value = "11/12/2023 23:12"
$ws.Cells.Item(69, 16).Value = 4
$ws.Cells.Item(69, 17).Value = "13/12/2023 10:58"
$ws.Cells.Item(69, 18).Value = 1.47
$ws.Cells.Item(69, 19).Value = "11/12/2023 23:12"
$ws.Cells.Item(69, 20).Value = 1.56
$ws.Cells.Item(69, 21).Value = "13/12/2023 10:51"
$ws.Cells.Item(69, 22).Value = "https://www.betexplorer.com/football/malta/premier-league/gudja-birkirkara/xp4RTBh6/"

# Row 70
$ws.Cells.Item(70, 1).Value = 69
$ws.Cells.Item(70, 2).Value = "malta"
$ws.Cells.Item(70, 3).Value = "premier-league"
$ws.Cells.Item(70, 4).Value = "2023-2024"
$ws.Cells.Item(70, 5).Value = 45273.58333333334
$ws.Cells.Item(70, 6).Value = "Hamrun"
$ws.Cells.Item(70, 7).Value = 2
$ws.Cells.Item(70, 8).Value = "Gzira"
$ws.Cells.Item(70, 9).Value = 1
$ws.Cells.Item(70, 10).Value = 1.75
$ws.Cells.Item(70, 11).Value = "12/12/2023 02:12"
$ws.Cells.Item(70, 12).Value = 1.56
$ws.Cells.Item(70, 13).Value = "13/12/2023 13:56"
$ws.Cells.Item(70, 14).Value = 3.31
$ws.Cells.Item(70, 15).Value = "12/12/2023 02:12"
$ws.Cells.Item(70, 16).Value = 3.74
$ws.Cells.Item(70, 17).Value = "13/12/2023 13:59"
$ws.Cells.Item(70, 18).Value = 4.29
$ws.Cells.Item(70, 19).Value = "12/12/2023 02:12"
$ws.Cells.Item(70, 20).Value = 6.23
$ws.Cells.Item(70, 21).Value = "13/12/2023 13:59"
$ws.Cells.Item(70, 22).Value = "https://www.betexplorer.com/football/malta/premier-league/hamrun-gzira/Gz3VSV7C/"

# Row 71
$ws.Cells.Item(71, 1).Value = 70
$ws.Cells.Item(71, 2).Value = "malta"
$ws.Cells.Item(71, 3).Value = "premier-league"
$ws.Cells.Item(71, 4).Value = "2023-2024"
$ws.Cells.Item(71, 5).Value = 45276.58333333334
$ws.Cells.Item(71, 6).Value = "Birkirkara"
$ws.Cells.Item(71, 7).Value = 3
$ws.Cells.Item(71, 8).Value = "Gzira"
$ws.Cells.Item(71, 9).Value = 2
$ws.Cells.Item(71, 10).Value = 2.21
$ws.Cells.Item(71, 11).Value = "15/12/2023 02:12"
$ws.Cells.Item(71, 12).Value = 2.03
$ws.Cells.Item(71, 13).Value = "16/12/2023 13:36"
$ws.Cells.Item(71, 14).Value = 3.01
$ws.Cells.Item(71, 15).Value = "15/12/2023 02:12"
$ws.Cells.Item(71, 16).Value = 3.2
$ws.Cells.Item(71, 17).Value = "16/12/2023 13:36"
$ws.Cells.Item(71, 18).Value = 3.12
$ws.Cells.Item(71, 19).Value = "15/12/2023 02:12"
$ws.Cells.Item(71, 20).Value = 3.76
$ws.Cells.Item(71, 21).Value = "16/12/2023 13:36"
$ws.Cells.Item(71, 22).Value = "https://www.betexplorer.com/football/malta/premier-league/birkirkara-gzira/MmYVQ9xP/"

# Row 72
$ws.Cells.Item(72, 1).Value = 71
$ws.Cells.Item(72, 2).Value = "malta"
$ws.Cells.Item(72, 3).Value = "premier-league"
$ws.Cells.Item(72, 4).Value = "2023-2024"
$ws.Cells.Item(72, 5).Value = 45276.58333333334
$ws.Cells.Item(72, 6).Value = "Gudja"
$ws.Cells.Item(72, 7).Value = 1
$ws.Cells.Item(72, 8).Value = "Sirens"
$ws.Cells.Item(72, 9).Value = 2
$ws.Cells.Item(72, 10).Value = 2.38
$ws.Cells.Item(72, 11).Value = "15/12/2023 02:12"
$ws.Cells.Item(72, 12).Value = 3.03
$ws.Cells.Item(72, 13).Value = "16/12/2023 13:50"
$ws.Cells.Item(72, 14).Value = 2.94
$ws.Cells.Item(72, 15).Value = "15/12/2023 02:12"
$ws.Cells.Item(72, 16).Value = 3.08
$ws.Cells.Item(72, 17).Value = "16/12/2023 13:50"
$ws.Cells.Item(72, 18).Value = 2.88
$ws.Cells.Item(72, 19).Value = "15/12/2023 02:12"
$ws.Cells.Item(72, 20).Value = 2.4
$ws.Cells.Item(72, 21).Value = "16/12/2023 13:50"
$ws.Cells.Item(72, 22).Value = "https://www.betexplorer.com/football/malta/premier-league/gudja-sirens/rcZRRkNI/"

# Row 73
$ws.Cells.Item(73, 1).Value = 72
$ws.Cells.Item(73, 2).Value = "malta"
$ws.Cells.Item(73, 3).Value = "premier-league"
$ws.Cells.Item(73, 4).Value = "2023-2024"
$ws.Cells.Item(73, 5).Value = 45276.67708333334
$ws.Cells.Item(73, 6).Value = "Balzan"
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = "Valletta"
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 2.52
$ws.Cells.Item(73, 11).Value = "15/12/2023 04:42"
$ws.Cells.Item(73, 12).Value = 3.01
$ws.Cells.Item(73, 13).Value = "16/12/2023 14:18"
$ws.Cells.Item(73, 14).Value = 2.9
$ws.Cells.Item(73, 15).Value = "15/12/2023 04:42"
$ws.Cells.Item(73, 16).Value = 3.18
$ws.Cells.Item(73, 17).Value = "16/12/2023 14:18"
$ws.Cells.Item(73, 18).Value = 2.73
$ws.Cells.Item(73, 19).Value = "15/12/2023 04:42"
$ws.Cells.Item(73, 20).Value = 2.35
$ws.Cells.Item(73, 21).Value = "16/12/2023 14:18"
$ws.Cells.Item(73, 22).Value = "https://www.betexplorer.com/football/malta/premier-league/balzan-fc-valletta/ziMDkgFP/"

# Row 74
$ws.Cells.Item(74, 1).Value = 73
$ws.Cells.Item(74, 2).Value = "malta"
$ws.Cells.Item(74, 3).Value = "premier-league"
$ws.Cells.Item(74, 4).Value = "2023-2024"
$ws.Cells.Item(74, 5).Value = 45277.45833333334
$ws.Cells.Item(74, 6).Value = "Naxxar"
$ws.Cells.Item(74, 7).Value = 1
$ws.Cells.Item(74, 8).Value = "Sliema"
$ws.Cells.Item(74, 9).Value = 4
$ws.Cells.Item(74, 10).Value = 3.58
$ws.Cells.Item(74, 11).Value = "15/12/2023 23:12"
$ws.Cells.Item(74, 12).Value = 4.74
$ws.Cells.Item(74, 13).Value = "17/12/2023 10:57"
$ws.Cells.Item(74, 14).Value = 3.16
$ws.Cells.Item(74, 15).Value = "15/12/2023 23:12"
$ws.Cells.Item(74, 16).Value = 3.5
$ws.Cells.Item(74, 17).Value = "17/12/2023 10:57"
$ws.Cells.Item(74, 18).Value = 1.96
$ws.Cells.Item(74, 19).Value = "15/12/2023 23:12"
$ws.Cells.Item(74, 20).Value = 1.74
$ws.Cells.Item(74, 21).Value = "17/12/2023 10:57"
$ws.Cells.Item(74, 22).Value = "https://www.betexplorer.com/football/malta/premier-league/naxxar-lions-sliema/CCDuoDNt/"

# Row 75
$ws.Cells.Item(75, 1).Value = 74
$ws.Cells.Item(75, 2).Value = "malta"
$ws.Cells.Item(75, 3).Value = "premier-league"
$ws.Cells.Item(75, 4).Value = "2023-2024"
$ws.Cells.Item(75, 5).Value = 45277.58333333334
$ws.Cells.Item(75, 6).Value = "Santa Lucia"
$ws.Cells.Item(75, 7).Value = 1
$ws.Cells.Item(75, 8).Value = "Hamrun"
$ws.Cells.Item(75, 9).Value = 3
$ws.Cells.Item(75, 10).Value = 10.97
$ws.Cells.Item(75, 11).Value = "16/12/2023 02:12"
$ws.Cells.Item(75, 12).Value = 8.37
$ws.Cells.Item(75, 13).Value = "17/12/2023 13:55"
$ws.Cells.Item(75, 14).Value = 5.37
$ws.Cells.Item(75, 15).Value = "16/12/2023 02:12"
$ws.Cells.Item(75, 16).Value = 4.67
$ws.Cells.Item(75, 17).Value = "17/12/2023 13:55"
$ws.Cells.Item(75, 18).Value = 1.2
$ws.Cells.Item(75, 19).Value = "16/12/2023 02:12"
$ws.Cells.Item(75, 20).Value = 1.36
$ws.Cells.Item(75, 21).Value = "17/12/2023 13:55"
$ws.Cells.Item(75, 22).Value = "https://www.betexplorer.com/football/malta/premier-league/santa-lucia-hamrun/YR1lqihg/"

# Row 76
$ws.Cells.Item(76, 1).Value = 75
$ws.Cells.Item(76, 2).Value = "malta"
$ws.Cells.Item(76, 3).Value = "premier-league"
$ws.Cells.Item(76, 4).Value = "2023-2024"
$ws.Cells.Item(76, 5).Value = 45277.58333333334
$ws.Cells.Item(76, 6).Value = "Hibernians"
$ws.Cells.Item(76, 7).Value = 0
$ws.Cells.Item(76, 8).Value = "Mosta"
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 10).Value = 1.75
$ws.Cells.Item(76, 11).Value = "16/12/2023 02:12"
$ws.Cells.Item(76, 12).Value = 1.65
$ws.Cells.Item(76, 13).Value = "17/12/2023 13:09"
$ws.Cells.Item(76, 14).Value = 3.41
$ws.Cells.Item(76, 15).Value = "16/12/2023 02:12"
$ws.Cells.Item(76, 16).Value = 3.63
$ws.Cells.Item(76, 17).Value = "17/12/2023 13:50"
$ws.Cells.Item(76, 18).Value = 4.13
$ws.Cells.Item(76, 19).Value = "16/12/2023 02:12"
$ws.Cells.Item(76, 20).Value = 5.22
$ws.Cells.Item(76, 21).Value = "17/12/2023 13:50"
$ws.Cells.Item(76, 22).Value = "https://www.betexplorer.com/football/malta/premier-league/hibernians-mosta-fc/fy3ppXxm/"

# Row 77
$ws.Cells.Item(77, 1).Value = 76
$ws.Cells.Item(77, 2).Value = "malta"
$ws.Cells.Item(77, 3).Value = "premier-league"
$ws.Cells.Item(77, 4).Value = "2023-2024"
$ws.Cells.Item(77, 5).Value = 45277.67708333334
$ws.Cells.Item(77, 6).Value = "Floriana"
$ws.Cells.Item(77, 7).Value = 3
$ws.Cells.Item(77, 8).Value = "Marsaxlokk"
$ws.Cells.Item(77, 9).Value = 2
$ws.Cells.Item(77, 10).Value = 1.83
$ws.Cells.Item(77, 11).Value = "16/12/2023 03:42"
$ws.Cells.Item(77, 12).Value = 1.63
$ws.Cells.Item(77, 13).Value = "17/12/2023 16:13"
$ws.Cells.Item(77, 14).Value = 3.21
$ws.Cells.Item(77, 15).Value = "16/12/2023 03:42"
$ws.Cells.Item(77, 16).Value = 3.51
$ws.Cells.Item(77, 17).Value = "17/12/2023 16:13"
$ws.Cells.Item(77, 18).Value = 4.03
$ws.Cells.Item(77, 19).Value = "16/12/2023 03:42"
$ws.Cells.Item(77, 20).Value = 5.76
$ws.Cells.Item(77, 21).Value = "17/12/2023 16:13"
$ws.Cells.Item(77, 22).Value = "https://www.betexplorer.com/football/malta/premier-league/floriana-marsaxlokk/216hrB7a/"

# --- Copy formatting for new rows from row 59 (A and E columns) ---
$ws.Range("A59").Copy() | Out-Null
$ws.Range("A60:A77").PasteSpecial(-4122) | Out-Null
$ws.Range("E59").Copy() | Out-Null
$ws.Range("E60:E77").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Update dimension reference ---
$ws.Range("E60:E77").NumberFormat = "YYYY-MM-DD HH:MM:SS"
